$wb = $excel.ActiveWorkbook

# --- Adform sheet: no longer the active tab, two rows' "code"/id values refreshed ---
$ws1 = $wb.Worksheets.Item("Adform")
$ws1.Range("A3").ClearContents()
$ws1.Range("B3").Value = 20181130001
$ws1.Range("A4").ClearContents()
$ws1.Range("B4").Value = 20181130001
$ws1.Range("B4").Select()

# --- AppNexus sheet: becomes the active tab, segment data updated for new test segments ---
$ws3 = $wb.Worksheets.Item("AppNexus")
$ws3.Range("A3").Value = 15662758
$ws3.Range("B3").Value = -2018113011
$ws3.Range("C3").Value = "Test Segment Name 30 Nov 2018 1"
$ws3.Range("A4").Value = 15662759
$ws3.Range("B4").Value = -2018113012
$ws3.Range("C4").Value = "Test Segment Name 30 Nov 2018 2"
$ws3.Range("H4").Value = $true

$ws3.Activate()
$ws3.Range("B5").Select()
